# Add three new submission rows (68-70) to the first worksheet
# ("八位序列号收集收集结果yd5" / sheet1.xml), mirroring the rows that were
# appended upstream. Columns: A=submitter, B=submit time, C=serial, D=QQ.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 67

$rows = @(
    @{ A = "如果"; B = 45949.5079398148; C = "f3bb7437"; D = "2752741520" },
    @{ A = "Yuki"; B = 45949.5953240741; C = "d492fa67"; D = "571584956" },
    @{ A = "隽昊"; B = 45949.7996180556; C = "57c84d5."; D = "992333681" }
)

# Scratch cell used to stage numeric-looking text so it lands in the sheet
# as a genuine string (shared-string) cell - matching how the existing QQ
# numbers in column D are stored - instead of being auto-coerced into a
# number by a plain `.Value =` assignment.
$scratch = $ws.Cells.Item(500, 10)

function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = "=""" + $escaped + """"
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

$r = $lastRow
foreach ($row in $rows) {
    $r = $r + 1

    $ws.Cells.Item($r, 1).Value = $row.A

    $dateCell = $ws.Cells.Item($r, 2)
    $ws.Cells.Item($lastRow, 2).Copy()
    $dateCell.PasteSpecial(-4122)
    $dateCell.Value = $row.B

    Set-TextValue $ws.Cells.Item($r, 3) $row.C
    Set-TextValue $ws.Cells.Item($r, 4) $row.D
}

$scratch.Clear()

Write-Output "done"
